# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps for the
# 67741593-67a5-440a-95f7-1285f266c3b8 row (row 5) on both the
# zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D5").Value = "2016-02-25 09:08:27"
$zhcn.Range("G5").Value = "2016-02-25 09:09:10"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D5").Value = "2016-02-25 09:08:37"
$dede.Range("G5").Value = "2016-02-25 09:09:27"
